# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the aggregated "全部类型" sheet, matching the regenerated site output.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetExhibition.Range("F2").Value = 586
$sheetExhibition.Range("F3").Value = 126
$sheetExhibition.Range("F7").Value = 1627

$sheetAll = $wb.Worksheets.Item("全部类型")
$sheetAll.Range("F2").Value = 586
$sheetAll.Range("F3").Value = 126
$sheetAll.Range("F11").Value = 1627
